$d = $word.ActiveDocument

# 1. Rename the UC007 use case from "Inserir Setor" to
#    "Enviar Mensagem Lembrete Compromisso".
$d.Content.Find.Execute("Inserir Setor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Enviar Mensagem Lembrete Compromisso", 2)

# 2. Strike through the UC008 row ("UC008" / "Manter Setor") to mark it as
#    removed/obsolete, applying the formatting to both the cell text and the
#    paragraph mark.
$tbl = $d.Tables.Item(1)
$tbl.Cell(9, 1).Range.Font.StrikeThrough = $true
$tbl.Cell(9, 2).Range.Font.StrikeThrough = $true
